# Auto-generated edit script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9984.75
$ws.Range("J32").Value = 10875
$ws.Range("L32").Value = 10875
$ws.Range("N32").Value = -11527
$ws.Range("H33").Value = 11364248
$ws.Range("I33").Value = 19231228
$ws.Range("J33").Value = 831.55554
$ws.Range("K33").Value = 19231228
$ws.Range("L33").Value = 831.55554
$ws.Range("M33").Value = -19230999
$ws.Range("N33").Value = -1289.55554
$ws.Range("H51").Value = 5199.75
$ws.Range("I51").Value = 4624
$ws.Range("J51").Value = 5391.6665
$ws.Range("K51").Value = 4624
$ws.Range("L51").Value = 5391.6665
$ws.Range("M51").Value = -4140
$ws.Range("N51").Value = -6359.6665
$ws.Range("H76").Value = 9455.77
$ws.Range("J76").Value = 10844.5
$ws.Range("L76").Value = 10844.5
$ws.Range("N76").Value = -11474.5
$ws.Range("H79").Value = 9455.77
$ws.Range("J79").Value = 10844.5
$ws.Range("L79").Value = 10844.5
$ws.Range("N79").Value = -13028.5
$ws.Range("H86").Value = 6247.769
$ws.Range("I86").Value = 6633.727
$ws.Range("J86").Value = 4125
$ws.Range("K86").Value = 6633.727
$ws.Range("L86").Value = 4125
$ws.Range("M86").Value = -5510.727
$ws.Range("N86").Value = -6371
$ws.Range("H89").Value = 6247.769
$ws.Range("I89").Value = 6633.727
$ws.Range("J89").Value = 4125
$ws.Range("K89").Value = 33168.635
$ws.Range("L89").Value = 20625
$ws.Range("M89").Value = -27552.635
$ws.Range("N89").Value = -31857
$ws.Range("H137").Value = 7777.841
$ws.Range("I137").Value = 3659.476
$ws.Range("J137").Value = 11538.087
$ws.Range("K137").Value = 10978.428
$ws.Range("L137").Value = 34614.261
$ws.Range("M137").Value = -8428.428
$ws.Range("N137").Value = -39714.261
$ws.Range("H138").Value = 3670.6667
$ws.Range("I138").Value = 3417.5334
$ws.Range("K138").Value = 10252.6002
$ws.Range("M138").Value = -5112.600199999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 91.5
$ws.Range("I4").Value = 91.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 91.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 24.5
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 334.33334
$ws.Range("I5").Value = 350.5
$ws.Range("J5").Value = 302
$ws.Range("K5").Value = 350.5
$ws.Range("L5").Value = 302
$ws.Range("M5").Value = -238.5
$ws.Range("N5").Value = -526
$ws.Range("H9").Value = 424242
$ws.Range("J9").Value = 424242
$ws.Range("L9").Value = 424242
$ws.Range("N9").Value = -424582
$ws.Range("H20").Value = 424242
$ws.Range("J20").Value = 424242
$ws.Range("L20").Value = 424242
$ws.Range("N20").Value = -424782
$ws.Range("H32").Value = 7429.516
$ws.Range("I32").Value = 4144.643
$ws.Range("K32").Value = 4144.643
$ws.Range("M32").Value = -3857.643
$ws.Range("H74").Value = 7950.028
$ws.Range("I74").Value = 2502.8386
$ws.Range("K74").Value = 2502.8386
$ws.Range("M74").Value = -1628.8386
$ws.Range("H77").Value = 7950.028
$ws.Range("I77").Value = 2502.8386
$ws.Range("K77").Value = 12514.193
$ws.Range("M77").Value = -8146.192999999999
$ws.Range("H110").Value = 4805.5835
$ws.Range("I110").Value = 5180.1816
$ws.Range("J110").Value = 685
$ws.Range("K110").Value = 5180.1816
$ws.Range("L110").Value = 685
$ws.Range("M110").Value = -3135.1816
$ws.Range("N110").Value = -4775
$ws.Range("H135").Value = 153870.67
$ws.Range("J135").Value = 153870.67
$ws.Range("L135").Value = 153870.67
$ws.Range("N135").Value = -164010.67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 334.33334
$ws.Range("I4").Value = 350.5
$ws.Range("J4").Value = 302
$ws.Range("K4").Value = 350.5
$ws.Range("L4").Value = 302
$ws.Range("M4").Value = -235.5
$ws.Range("N4").Value = -532
$ws.Range("H80").Value = 900.63635
$ws.Range("J80").Value = 886.9231
$ws.Range("L80").Value = 886.9231
$ws.Range("N80").Value = -2882.9231
$ws.Range("H83").Value = 900.63635
$ws.Range("J83").Value = 886.9231
$ws.Range("L83").Value = 4434.6155
$ws.Range("N83").Value = -14418.6155
$ws.Range("H94").Value = 2249
$ws.Range("I94").Value = 2249
$ws.Range("K94").Value = 2249
$ws.Range("M94").Value = -1798
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1142.875
$ws.Range("I7").Value = 877.6429000000001
$ws.Range("K7").Value = 877.6429000000001
$ws.Range("M7").Value = -764.6429000000001
$ws.Range("H22").Value = 1340.7368
$ws.Range("I22").Value = 703.875
$ws.Range("J22").Value = 1803.909
$ws.Range("K22").Value = 703.875
$ws.Range("L22").Value = 1803.909
$ws.Range("M22").Value = -353.875
$ws.Range("N22").Value = -2503.909
$ws.Range("H31").Value = 122964.9
$ws.Range("I31").Value = 202482.2
$ws.Range("J31").Value = 43447.6
$ws.Range("K31").Value = 202482.2
$ws.Range("L31").Value = 43447.6
$ws.Range("M31").Value = -202187.2
$ws.Range("N31").Value = -44037.6
$ws.Range("H34").Value = 122964.9
$ws.Range("I34").Value = 202482.2
$ws.Range("J34").Value = 43447.6
$ws.Range("K34").Value = 202482.2
$ws.Range("L34").Value = 43447.6
$ws.Range("M34").Value = -202280.2
$ws.Range("N34").Value = -43851.6
$ws.Range("H58").Value = 8980.118
$ws.Range("J58").Value = 17837
$ws.Range("L58").Value = 17837
$ws.Range("N58").Value = -18243
$ws.Range("H99").Value = 3352.3333
$ws.Range("I99").Value = 3334.125
$ws.Range("J99").Value = 3498
$ws.Range("K99").Value = 3334.125
$ws.Range("L99").Value = 3498
$ws.Range("M99").Value = -1836.125
$ws.Range("N99").Value = -6494
$ws.Range("H126").Value = 3352.3333
$ws.Range("I126").Value = 3334.125
$ws.Range("J126").Value = 3498
$ws.Range("K126").Value = 10002.375
$ws.Range("L126").Value = 10494
$ws.Range("M126").Value = -7532.375
$ws.Range("N126").Value = -15434
$ws.Range("H132").Value = 1409772.5
$ws.Range("I132").Value = 2674.8667
$ws.Range("K132").Value = 8024.6001
$ws.Range("M132").Value = -5494.6001
$ws.Range("H134").Value = 6529.8076
$ws.Range("I134").Value = 1893.0454
$ws.Range("K134").Value = 5679.1362
$ws.Range("M134").Value = -3144.1362
$ws.Range("H136").Value = 8980.118
$ws.Range("J136").Value = 17837
$ws.Range("L136").Value = 53511
$ws.Range("N136").Value = -58611
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 229.58333
$ws.Range("I7").Value = 185
$ws.Range("J7").Value = 292
$ws.Range("K7").Value = 555
$ws.Range("L7").Value = 876
$ws.Range("M7").Value = -443
$ws.Range("N7").Value = -1100
$ws.Range("H103").Value = 1637
$ws.Range("I103").Value = 1110
$ws.Range("J103").Value = 2164
$ws.Range("K103").Value = 3330
$ws.Range("L103").Value = 6492
$ws.Range("M103").Value = -2451
$ws.Range("N103").Value = -8250
$ws.Range("H131").Value = 1435.89
$ws.Range("I131").Value = 807.25
$ws.Range("J131").Value = 1462.0834
$ws.Range("K131").Value = 2421.75
$ws.Range("L131").Value = 4386.2502
$ws.Range("M131").Value = 2618.25
$ws.Range("N131").Value = -14466.2502
$ws.Range("H137").Value = 12611
$ws.Range("J137").Value = 15481.444
$ws.Range("L137").Value = 46444.33199999999
$ws.Range("N137").Value = -56644.33199999999
$ws.Range("H140").Value = 1982.8125
$ws.Range("I140").Value = 1463.4615
$ws.Range("K140").Value = 4390.3845
$ws.Range("M140").Value = 789.6154999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 64.333336
$ws.Range("I2").Value = 60.333332
$ws.Range("J2").Value = 68.333336
$ws.Range("K2").Value = 60.333332
$ws.Range("L2").Value = 68.333336
$ws.Range("M2").Value = 52.666668
$ws.Range("N2").Value = -294.333336
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4004.25
$ws.Range("I40").Value = 3346.2856
$ws.Range("K40").Value = 3346.2856
$ws.Range("M40").Value = -3210.2856
$ws.Range("H46").Value = 3690.1428
$ws.Range("I46").Value = 3292.8572
$ws.Range("J46").Value = 3888.7856
$ws.Range("K46").Value = 3292.8572
$ws.Range("L46").Value = 3888.7856
$ws.Range("M46").Value = -3104.8572
$ws.Range("N46").Value = -4264.7856
$ws.Range("H80").Value = 16666.666
$ws.Range("H83").Value = 16666.666
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 1248867.4
$ws.Range("I132").Value = 2528.875
$ws.Range("K132").Value = 7586.625
$ws.Range("M132").Value = -5056.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 6499.4443
$ws.Range("I13").Value = 6498
$ws.Range("J13").Value = 6499.625
$ws.Range("K13").Value = 6498
$ws.Range("L13").Value = 6499.625
$ws.Range("M13").Value = -6358
$ws.Range("N13").Value = -6779.625
$ws.Range("H136").Value = 237847.98
$ws.Range("I136").Value = 1410.4348
$ws.Range("J136").Value = 485032.7
$ws.Range("K136").Value = 4231.3044
$ws.Range("L136").Value = 1455098.1
$ws.Range("M136").Value = -1681.3044
$ws.Range("N136").Value = -1460198.1
